$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.67%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'37.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.56%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.126"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.23%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07777"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.57%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-1.31%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.224"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.61%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.885"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-8.44%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-11.69%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9218"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.22%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1217"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.99%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1912"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.54%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09217"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.99%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'-0.58%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09683"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.37%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001366"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.74%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005954"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-6.31%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.559"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.52%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3406"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.79%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'5.258"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.11%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1268"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.67%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2591"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.58%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'5,590.67%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04369"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.21%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-3.30%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-8.69%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02090"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-5.21%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05071"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.36%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007692"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.99%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009815"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.63%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1347"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.89%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'1.60%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009585"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'8.11%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006681"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.32%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.91%"
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'BOLO"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'0.002937"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-2.92%"
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'CoinbaseStockToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.001201"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.49%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.91%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.91%"
$ws.Range("E51").Style = "Normal"
